# Apply sprint backlog updates to "Sprint 2" sheet after scrum.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2")

# New backlog stories added after the scrum (rows 23-28).
# Values are entered in the same order the author typed them so the
# shared-string table ends up in the same sequence.
$ws.Cells.Item(23,1).Value = "56, Create user interface for creating a schedule, adding classes, and saving a schedule"
$ws.Cells.Item(23,2).Value = "1,2"
$ws.Cells.Item(23,3).Value = "F"
$ws.Cells.Item(23,4).Value = "T"
$ws.Cells.Item(23,5).HorizontalAlignment = -4152

$ws.Cells.Item(24,1).Value = "57, Create user interface for recalling schedules, schedule selection, and detail display"
$ws.Cells.Item(24,2).Value = "2,3"
$ws.Cells.Item(24,3).Value = "F"
$ws.Cells.Item(24,4).Value = "T"

$ws.Cells.Item(25,1).Value = "58, Implement backend functionality for saving and recalling saved schedules"
$ws.Cells.Item(25,2).Value = "2,3"
$ws.Cells.Item(25,3).Value = "F"
$ws.Cells.Item(25,4).Value = "T"

$ws.Cells.Item(26,1).Value = "59, Implement data structure for mapping building prefixes to physical addresses"
$ws.Cells.Item(26,3).Value = "F"
$ws.Cells.Item(26,4).Value = "T"

$ws.Cells.Item(27,1).Value = "60, Implement driver code to utilize google maps API to generate route from provided schedule detail"
$ws.Cells.Item(27,3).Value = "F"
$ws.Cells.Item(27,4).Value = "T"

$ws.Cells.Item(26,2).Value = "1,5"
$ws.Cells.Item(27,2).Value = "1,5"

$ws.Cells.Item(28,1).Value = "61, Implement front end for route display after API calls to generate route"
$ws.Cells.Item(28,2).Value = "1,5"
$ws.Cells.Item(28,3).Value = "F"
$ws.Cells.Item(28,4).Value = "T"

# Widen column A to fit the longer story text and resize the view.
$ws.Columns.Item(1).ColumnWidth = 110

# Scroll the view back up and leave the selection on the last edited cell,
# matching where the author left off.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D28").Select()
